$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F4 (想去人数 for 南宁·草莓动漫节) and F6 (南宁·布谷鸟动漫展4th)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1635
$wsExhibit.Range("F6").Value = 56

# Sheet "全部类型" (sheet4): same two events appear at F4 and F7
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1635
$wsAll.Range("F7").Value = 56
